$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 53, shifting existing rows 53:62 down to 54:63.
$ws.Rows("53:53").Insert()

# Populate the new row 53 with the same constant columns used by the
# surrounding "Camote" / "Vega Modelo de Temuco" records, plus the new
# weekly data point.
$ws.Range("A53").Value = 10
$ws.Range("B53").Value = "Vega Modelo de Temuco"
$ws.Range("C53").Value = "La Araucanía"
$ws.Range("D53").Value = 44642
$ws.Range("E53").Value = 9
$ws.Range("F53").Value = 100114002
$ws.Range("G53").Value = "Camote"
$ws.Range("H53").Value = "Sin especificar"
$ws.Range("I53").Value = "Primera"
$ws.Range("J53").Value = 25
$ws.Range("K53").Value = 18000
$ws.Range("L53").Value = 18000
$ws.Range("M53").Value = 18000
$ws.Range("N53").Value = "$/malla 20 kilos"
$ws.Range("O53").Value = "Perú"
$ws.Range("P53").Value = 900
$ws.Range("Q53").Value = 20
$ws.Range("R53").Value = "Hortaliza"
